# Generate Report for Handback
# Refresh the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on row 2 of the
# per-language report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 07:20:24"
$wsZhCn.Range("H2").Value = "2016-03-24 07:20:47"

$wsDeDe = $wb.Worksheets("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 07:20:28"
$wsDeDe.Range("H2").Value = "2016-03-24 07:20:54"
